{"js": "// Commit: \"Edited Witch and witch prototype fold\"\n//\n// Change: \"...individuals would contribute to the system more smoothly. \"\n//      -> \"...individuals would contribute to the capital system more smoothly.\"\n//\n// i.e. the word \"capital \" is inserted before \"system more smoothly.\" and the\n// sentence's former trailing space (just before the bookmark) is dropped.\n\nconst body = context.document.body;\n\n// Search for the exact phrase that changes, trailing space included, so the\n// replacement can also absorb/drop that trailing space in one shot.\nconst searchText = \"the system more smoothly. \";\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Target text not found: ${JSON.stringify(searchText)}`);\n}\n\nconst target = results.items[0];\ntarget.insertText(\"the capital system more smoothly.\", \"Replace\");\nawait context.sync();\n", "ps1": "# Commit: \"Edited Witch and witch prototype fold\"\n#\n# Change: \"...individuals would contribute to the system more smoothly. \"\n#      -> \"...individuals would contribute to the capital system more smoothly.\"\n#\n# i.e. the word \"capital \" is inserted before \"system more smoothly.\" and the\n# sentence's former trailing space (just before the bookmark) is dropped.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"the system more smoothly. \"\n$find.Replacement.Text = \"the capital system more smoothly.\"\n\n$find.Execute(\n    $find.Text,       # FindText\n    $false,           # MatchCase\n    $false,           # MatchWholeWord\n    $false,           # MatchWildcards\n    $false,           # MatchSoundsLike\n    $false,           # MatchAllWordForms\n    $true,            # Forward\n    $wdFindContinue,  # Wrap\n    $false,           # Format\n    $find.Replacement.Text,  # ReplaceWith\n    $wdReplaceAll     # Replace\n)\n"}
